# Weekly "Fruta / hortaliza" update for the Cilantro sheet.
#
# Two new daily price rows (market date 2021-10-07 = serial 44476) are
# inserted at the top of the historical data block (row 287), pushing the
# existing 287-348 rows down to 289-350 (dimension grows from R348 to R350).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 287; this shifts every
# existing row (287..348) down by two (-> 289..350) and carries the
# date-column (D) number format down onto the freshly inserted rows,
# exactly like Excel's native "Insert Copied/Blank Rows" behaviour.
$ws.Rows("287:288").Insert()

# Row 287 - "$/caja 36 atados" (box of 36 bunches) record for 2021-10-07
$ws.Range("A287").Value = 9
$ws.Range("B287").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C287").Value = "Metropolitana"
$ws.Range("D287").Value = 44476
$ws.Range("E287").Value = 13
$ws.Range("F287").Value = 100112040
$ws.Range("G287").Value = "Cilantro"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Primera"
$ws.Range("J287").Value = 52
$ws.Range("K287").Value = 4000
$ws.Range("L287").Value = 4000
$ws.Range("M287").Value = 4000
$ws.Range("N287").Value = "`$/caja 36 atados"
$ws.Range("O287").Value = "Región Metropolitana"
$ws.Range("P287").Value = 111
$ws.Range("Q287").Value = 36
$ws.Range("R287").Value = "Hortaliza"

# Row 288 - "$/docena de atados" (dozen of bunches) record for 2021-10-07
$ws.Range("A288").Value = 9
$ws.Range("B288").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C288").Value = "Metropolitana"
$ws.Range("D288").Value = 44476
$ws.Range("E288").Value = 13
$ws.Range("F288").Value = 100112040
$ws.Range("G288").Value = "Cilantro"
$ws.Range("H288").Value = "Sin especificar"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 133
$ws.Range("K288").Value = 7000
$ws.Range("L288").Value = 8000
$ws.Range("M288").Value = 7496
$ws.Range("N288").Value = "`$/docena de atados"
$ws.Range("O288").Value = "Región Metropolitana"
$ws.Range("P288").Value = 2499
$ws.Range("Q288").Value = 3
$ws.Range("R288").Value = "Hortaliza"
